$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 gets the data that used to be in row 3, row 3 gets the data that used to be in row 2
# (columns A-H swap their contents between the two rows)
$ws.Range("A2").Value = 111661765
$ws.Range("B2").Value = 89953
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value = 3884
$ws.Range("F2").Value = "Hasselticka"
$ws.Range("G2").Value = "Dichomitus campestris"
$ws.Range("H2").Value = "(Quél.) Domański & Orlicz"

$ws.Range("A3").Value = 111661750
$ws.Range("B3").Value = 99581
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 221317
$ws.Range("F3").Value = "Gullklöver"
$ws.Range("G3").Value = "Trifolium aureum"
$ws.Range("H3").Value = "Pollich"

# Coordinates in Q/R get rounded to whole numbers for both rows
$ws.Range("Q2").Value = 686742
$ws.Range("R2").Value = 6619854
$ws.Range("Q3").Value = 686742
$ws.Range("R3").Value = 6619854

# Time columns (Z, AB) are cleared out for both rows
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()

# Biotop (AH) moves from row 2 to row 3
$ws.Range("AH2").ClearContents()
$ws.Range("AH3").Value = "Vägkant"
